# Wrap every remaining blank-template placeholder (and the two dates) in
# square brackets, e.g. MMMM -> [MMMM], so the template makes it obvious
# which spans still need to be filled in by hand.

$d = $word.ActiveDocument
$tab = [char]9

# Simple four-letter (ish) placeholder tokens that just need to become
# [TOKEN] wherever they occur. MatchWholeWord avoids any chance of a
# partial-word hit, MatchCase keeps us from touching ordinary prose.
$placeholders = @(
    "MMMM", "NNNN", "OOOO", "PPPP", "QQQQ", "RRRR", "SSSS", "TTTT", "UUUU",
    "GGGG", "XXXY", "XXYY", "IIII", "JJJJ", "KKKK", "LLLL"
)

foreach ($token in $placeholders) {
    $d.Content.Find.Execute($token, $true, $true, $false, $false, $false, `
        $true, 1, $false, "[$token]", 2) | Out-Null
}

# "Registro <tab>Profissional AAAA" collapses the stray double tab into a
# single one and wraps AAAA at the same time.
$d.Content.Find.Execute( `
    "${tab}Registro ${tab}Profissional AAAA", $true, $false, $false, $false, `
    $false, $true, 1, $false, "${tab}Registro Profissional [AAAA]", 2) | Out-Null

# The "inerentes a" split (left over from an old spell-check flag on the
# lone word "a") is re-typed as a single run of plain text, which also
# clears the now-stale proofErr markers around it.
$d.Content.Find.Execute( `
    "poderes, inerentes a o bom e fiel cumprimento deste mandato, podendo representar o OUTORGANTE perante a [LLLL],", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "poderes, inerentes a o bom e fiel cumprimento deste mandato, podendo representar o OUTORGANTE perante a [LLLL],", `
    2) | Out-Null

# The two dates at the bottom of the document get bracketed as a whole
# (including the embedded " / " and "de ... de" text), not just the digits.
$d.Content.Find.Execute( `
    "Esta procuração é válida até 14 / 07 / 2026.", $true, $false, $false, `
    $false, $false, $true, 1, $false, `
    "Esta procuração é válida até [14 / 07 / 2026].", 2) | Out-Null

$d.Content.Find.Execute( `
    "Boa Vista, RR 14 de julho de 2025.", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Boa Vista, RR [14 de julho de 2025].", 2) | Out-Null
